$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 112, shifting the existing rows 112:186 down to 113:187
# (mirrors Excel's Rows("112").Insert / xlShiftDown behaviour)
$ws.Rows("112").Insert()

# Populate the newly inserted row 112 with the new data record
$ws.Range("A112").Value = 10
$ws.Range("B112").Value = "Vega Modelo de Temuco"
$ws.Range("C112").Value = "La Araucanía"
$ws.Range("D112").Value = 44438
$ws.Range("E112").Value = 9
$ws.Range("F112").Value = 100112037
$ws.Range("G112").Value = "Cebollín"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 30
$ws.Range("K112").Value = 9000
$ws.Range("L112").Value = 9000
$ws.Range("M112").Value = 9000
$ws.Range("N112").Value = "$/docena de paquetes"
$ws.Range("O112").Value = "Provincia de Cautín"
$ws.Range("P112").Value = 750
$ws.Range("Q112").Value = 12
$ws.Range("R112").Value = "Hortaliza"
